$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error: Marking row (B11/C11) and Total row (B12/C12/E12)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -14
$ws.Range("E12").Value = "66 / 112"
